$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, shifting existing data rows down by one.
$ws.Rows("2:2").Insert()

# Populate the new row 2 with the new first data point (date 2020-03-05 = serial 43895).
$ws.Cells.Item(2, 1).Value = 43895
$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 0

# Carry over the date number format from the row below, matching the workbook's
# existing date style instead of minting a new one.
$ws.Cells.Item(3, 1).Copy()
$ws.Cells.Item(2, 1).PasteSpecial(-4122)

# Update the sheet view: reset the scrolled-to cell and change the active selection.
$ws.Range("D4").Select()
